# Apply updated 'want to go' (想去人数) counts scraped at commit 456a3b4
# Values only change on sheets 1 (展览), 2 (演出) and 4 (全部类型); sheet 3 (本地生活) is untouched.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 1138
$ws.Range("F3").Value = 1956
$ws.Range("F4").Value = 615
$ws.Range("F5").Value = 1266
$ws.Range("F7").Value = 46
$ws.Range("F9").Value = 341
$ws.Range("F10").Value = 123
$ws.Range("F11").Value = 103
$ws.Range("F12").Value = 836
$ws.Range("F13").Value = 256
$ws.Range("F14").Value = 132
$ws.Range("F17").Value = 346
$ws.Range("F18").Value = 247
$ws.Range("F19").Value = 707
$ws.Range("F20").Value = 80
$ws.Range("F21").Value = 671
$ws.Range("F22").Value = 202
$ws.Range("F24").Value = 912
$ws.Range("F25").Value = 367
$ws.Range("F26").Value = 197
$ws.Range("F28").Value = 307

$ws = $wb.Worksheets.Item(2)
$ws.Range("F12").Value = 26

$ws = $wb.Worksheets.Item(4)
$ws.Range("F3").Value = 1138
$ws.Range("F4").Value = 1956
$ws.Range("F5").Value = 615
$ws.Range("F6").Value = 1266
$ws.Range("F9").Value = 46
$ws.Range("F11").Value = 341
$ws.Range("F12").Value = 123
$ws.Range("F13").Value = 103
$ws.Range("F14").Value = 836
$ws.Range("F15").Value = 256
$ws.Range("F16").Value = 132
$ws.Range("F22").Value = 346
$ws.Range("F25").Value = 247
$ws.Range("F26").Value = 707
$ws.Range("F27").Value = 80
$ws.Range("F28").Value = 671
$ws.Range("F29").Value = 202
$ws.Range("F31").Value = 912
$ws.Range("F32").Value = 367
$ws.Range("F35").Value = 197
$ws.Range("F37").Value = 307
$ws.Range("F42").Value = 26

